$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue "D2" "242.65"
Set-TextValue "D3" "23.06"
Set-TextValue "D4" "5.427"
Set-TextValue "D5" "0.05895"
Set-TextValue "D6" "3.445"
Set-TextValue "D7" "6.524"
Set-TextValue "D8" "0.8100"
Set-TextValue "D9" "0.9686"
Set-TextValue "D10" "0.1413"
Set-TextValue "D11" "0.07427"
Set-TextValue "D12" "0.03274"
Set-TextValue "D13" "0.03052"
Set-TextValue "D14" "0.09338"
Set-TextValue "D15" "3.852"
Set-TextValue "D16" "0.001578"
Set-TextValue "D17" "0.04688"
Set-TextValue "D18" "0.0005911"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue "D19" "0.005876"
Set-TextValue "D20" "0.001253"
Set-TextValue "D21" "0.004896"
Set-TextValue "D22" "0.00006802"
Set-TextValue "D23" "3.592"
Set-TextValue "D25" "0.3221"
Set-TextValue "D27" "0.0002284"
Set-TextValue "D40" "0.03924"
Set-TextValue "D41" "0.006181"
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.009752"
Set-TextValue "D45" "0.00005193"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "D47" "0.6661"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue "D48" "0.002379"
